# Update leve profit figures across multiple sheets (scheduled market-price refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 695.3077
$ws.Range("I11").Value = 695.3077
$ws.Range("K11").Value = 695.3077
$ws.Range("M11").Value = -555.3077

$ws.Range("H113").Value = 3796.3157
$ws.Range("I113").Value = 2390.4546
$ws.Range("J113").Value = 5729.375
$ws.Range("K113").Value = 2390.4546
$ws.Range("L113").Value = 5729.375
$ws.Range("M113").Value = 863.5454
$ws.Range("N113").Value = -12237.375

$ws.Range("H135").Value = 527679.7
$ws.Range("I135").Value = 527679.7
$ws.Range("K135").Value = 4749117.3
$ws.Range("M135").Value = -4746582.3

$ws.Range("H137").Value = 1933.5526
$ws.Range("I137").Value = 1257.7587
$ws.Range("J137").Value = 4111.1113
$ws.Range("K137").Value = 3773.2761
$ws.Range("L137").Value = 12333.3339
$ws.Range("M137").Value = -1223.2761
$ws.Range("N137").Value = -17433.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2681.3125
$ws.Range("I2").Value = 1737.3636
$ws.Range("J2").Value = 4758
$ws.Range("K2").Value = 1737.3636
$ws.Range("L2").Value = 4758
$ws.Range("M2").Value = -1624.3636
$ws.Range("N2").Value = -4984

$ws.Range("H32").Value = 19570.09
$ws.Range("I32").Value = 4260.32
$ws.Range("J32").Value = 96118.92999999999
$ws.Range("K32").Value = 4260.32
$ws.Range("L32").Value = 96118.92999999999
$ws.Range("M32").Value = -3973.32
$ws.Range("N32").Value = -96692.92999999999

$ws.Range("H61").Value = 1383.9656
$ws.Range("I61").Value = 1328.7273
$ws.Range("J61").Value = 1557.5714
$ws.Range("K61").Value = 1328.7273
$ws.Range("L61").Value = 1557.5714
$ws.Range("M61").Value = -1116.7273
$ws.Range("N61").Value = -1981.5714

$ws.Range("H74").Value = 1113.5483
$ws.Range("I74").Value = 818.34784
$ws.Range("J74").Value = 1962.25
$ws.Range("K74").Value = 818.34784
$ws.Range("L74").Value = 1962.25
$ws.Range("M74").Value = 55.65215999999998
$ws.Range("N74").Value = -3710.25

$ws.Range("H77").Value = 1113.5483
$ws.Range("I77").Value = 818.34784
$ws.Range("J77").Value = 1962.25
$ws.Range("K77").Value = 4091.7392
$ws.Range("L77").Value = 9811.25
$ws.Range("M77").Value = 276.2608
$ws.Range("N77").Value = -18547.25

$ws.Range("H101").Value = 41867.332
$ws.Range("J101").Value = 41867.332
$ws.Range("L101").Value = 41867.332
$ws.Range("N101").Value = -48357.332

$ws.Range("H110").Value = 7665.4
$ws.Range("I110").Value = 7665.4
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 7665.4
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -5620.4
$ws.Range("N110").ClearContents()

$ws.Range("H116").Value = 2681.3125
$ws.Range("I116").Value = 1737.3636
$ws.Range("J116").Value = 4758
$ws.Range("K116").Value = 1737.3636
$ws.Range("L116").Value = 4758
$ws.Range("M116").Value = 556.6364000000001
$ws.Range("N116").Value = -9346

$ws.Range("H132").Value = 2463.6843
$ws.Range("I132").Value = 2175.75
$ws.Range("K132").Value = 6527.25
$ws.Range("M132").Value = -3997.25

$ws.Range("H136").Value = 1383.9656
$ws.Range("I136").Value = 1328.7273
$ws.Range("J136").Value = 1557.5714
$ws.Range("K136").Value = 3986.1819
$ws.Range("L136").Value = 4672.7142
$ws.Range("M136").Value = -1436.1819
$ws.Range("N136").Value = -9772.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2681.3125
$ws.Range("I3").Value = 1737.3636
$ws.Range("J3").Value = 4758
$ws.Range("K3").Value = 1737.3636
$ws.Range("L3").Value = 4758
$ws.Range("M3").Value = -1623.3636
$ws.Range("N3").Value = -4986

$ws.Range("H5").Value = 5550
$ws.Range("I5").Value = 800
$ws.Range("K5").Value = 800
$ws.Range("M5").Value = -687

$ws.Range("H53").Value = 40770
$ws.Range("J53").Value = 40770
$ws.Range("L53").Value = 40770
$ws.Range("N53").Value = -41918

$ws.Range("H134").Value = 1360.5238
$ws.Range("I134").Value = 1319.5264
$ws.Range("J134").Value = 1750
$ws.Range("K134").Value = 3958.5792
$ws.Range("L134").Value = 5250
$ws.Range("M134").Value = -1423.5792
$ws.Range("N134").Value = -10320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2504.6667
$ws.Range("I31").Value = 1657.3
$ws.Range("J31").Value = 3808.3076
$ws.Range("K31").Value = 1657.3
$ws.Range("L31").Value = 3808.3076
$ws.Range("M31").Value = -1362.3
$ws.Range("N31").Value = -4398.3076

$ws.Range("H34").Value = 2504.6667
$ws.Range("I34").Value = 1657.3
$ws.Range("J34").Value = 3808.3076
$ws.Range("K34").Value = 1657.3
$ws.Range("L34").Value = 3808.3076
$ws.Range("M34").Value = -1455.3
$ws.Range("N34").Value = -4212.3076

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 9818.75
$ws.Range("J93").Value = 9818.75
$ws.Range("L93").Value = 9818.75
$ws.Range("N93").Value = -13562.75

$ws.Range("H113").Value = 1615.75
$ws.Range("I113").Value = 1435.4667
$ws.Range("J113").Value = 2156.6
$ws.Range("K113").Value = 1435.4667
$ws.Range("L113").Value = 2156.6
$ws.Range("M113").Value = 734.5333000000001
$ws.Range("N113").Value = -6496.6

$ws.Range("H132").Value = 3452.8386
$ws.Range("I132").Value = 2629.0908
$ws.Range("J132").Value = 3905.9
$ws.Range("K132").Value = 7887.2724
$ws.Range("L132").Value = 11717.7
$ws.Range("M132").Value = -5357.2724
$ws.Range("N132").Value = -16777.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2878.5
$ws.Range("I61").Value = 2379.6667
$ws.Range("K61").Value = 2379.6667
$ws.Range("M61").Value = -2177.6667

$ws.Range("H113").Value = 2878.5
$ws.Range("I113").Value = 2379.6667
$ws.Range("K113").Value = 2379.6667
$ws.Range("M113").Value = -209.6667000000002

$ws.Range("H132").Value = 3330.513
$ws.Range("I132").Value = 3382.6667
$ws.Range("J132").Value = 3213.1667
$ws.Range("K132").Value = 10148.0001
$ws.Range("L132").Value = 9639.500100000001
$ws.Range("M132").Value = -7618.000100000001
$ws.Range("N132").Value = -14699.5001

$ws.Range("H136").Value = 2122.5417
$ws.Range("I136").Value = 1711.4762
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 5134.4286
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -2584.4286
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 20135.5
$ws.Range("I69").Value = 10000
$ws.Range("J69").Value = 30271
$ws.Range("K69").Value = 10000
$ws.Range("L69").Value = 30271
$ws.Range("M69").Value = -9251
$ws.Range("N69").Value = -31769

$ws.Range("H70").Value = 20000
$ws.Range("J70").Value = 20000
$ws.Range("L70").Value = 20000
$ws.Range("N70").Value = -20630

$ws.Range("H72").Value = 20135.5
$ws.Range("I72").Value = 10000
$ws.Range("J72").Value = 30271
$ws.Range("K72").Value = 30000
$ws.Range("L72").Value = 90813
$ws.Range("M72").Value = -26256
$ws.Range("N72").Value = -98301

$ws.Range("H73").Value = 20000
$ws.Range("J73").Value = 20000
$ws.Range("L73").Value = 20000
$ws.Range("N73").Value = -22184

$ws.Range("H113").Value = 31841.406
$ws.Range("I113").Value = 38679.73
$ws.Range("J113").Value = 2208.6667
$ws.Range("K113").Value = 116039.19
$ws.Range("L113").Value = 6626.000100000001
$ws.Range("M113").Value = -113869.19
$ws.Range("N113").Value = -10966.0001

$ws.Range("H136").Value = 1181.375
$ws.Range("I136").Value = 1126.4667
$ws.Range("J136").Value = 2005
$ws.Range("K136").Value = 3379.4001
$ws.Range("L136").Value = 6015
$ws.Range("M136").Value = -829.4000999999998
$ws.Range("N136").Value = -11115
